# corrette formule Msq per coda M/M/1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scroll the view so row 82 is the first visible row (topLeftCell="A82")
$ws.Application.ActiveWindow.ScrollRow = 82

# --- Block 1 (rows 26-33): change service-rate input B28 ---
$ws.Range("B28").Value = 3600

# --- Block 2 (rows 37-45): change arrival-rate input D38 ---
$ws.Range("D38").Value = 0.1

# --- Row heights for section headers at rows 73 and 85 ---
$ws.Rows("73").RowHeight = 17.25
$ws.Rows("85").RowHeight = 17.25

# --- Block 3 (rows 86-93): change service-rate input B88 ---
$ws.Range("B88").Value = 1200

# B90 no longer needed (formula removed) - clear the cell entirely
$ws.Range("B90").ClearContents()

# B92/B93 formulas simplified (M/M/1 queue correction)
$ws.Range("B92").Formula = "=B89*B88/(1-B89)"
$ws.Range("B93").Formula = "=B92*B86*D86"
